$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Saakshi"
$ws.Range("B4").Value = '$2b$12$XK/ATOEGuqgpLzTXeirhreOPRkFDaxwZJe5YnG9k.mo8CELGHTeLy'
